# Insert a new data row at row 264, shifting existing rows 264:371 down to 265:372,
# then populate the new row 264 with the updated record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 264 (this shifts row 264 -> 265, ..., row 371 -> 372)
$ws.Rows.Item(264).Insert()

# Populate the newly inserted row 264 with the new record.
$ws.Cells.Item(264, 1).Value = 3
$ws.Cells.Item(264, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(264, 3).Value = "Coquimbo"
$ws.Cells.Item(264, 4).Value = 44704
$ws.Cells.Item(264, 5).Value = 5
$ws.Cells.Item(264, 6).Value = 100112031
$ws.Cells.Item(264, 7).Value = "Poroto verde"
$ws.Cells.Item(264, 8).Value = "Magnum"
$ws.Cells.Item(264, 9).Value = "Primera"
$ws.Cells.Item(264, 10).Value = 73
$ws.Cells.Item(264, 11).Value = 27000
$ws.Cells.Item(264, 12).Value = 28000
$ws.Cells.Item(264, 13).Value = 27479
$ws.Cells.Item(264, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(264, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(264, 16).Value = 1099
$ws.Cells.Item(264, 17).Value = 25
$ws.Cells.Item(264, 18).Value = "Hortaliza"
